# Cap nhat toi uu
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Update existing A:D value cells (rows 3,4,6,7) ----
$ws.Range("A3").Value = 78
$ws.Range("B3").Value = 441
$ws.Range("C3").Value = 441
$ws.Range("D3").Value = 4921

$ws.Range("A4").Value = 403.65
$ws.Range("B4").Value = 1965.6
$ws.Range("C4").Value = 4068
$ws.Range("D4").Value = 16907.599999999999

$ws.Range("A6").Value = 51.52
$ws.Range("B6").Value = 59.21
$ws.Range("C6").Value = 65.010000000000005
$ws.Range("D6").Value = 66.849999999999994

$ws.Range("A7").Value = 49.16
$ws.Range("B7").Value = 63.32
$ws.Range("C7").Value = 68.22
$ws.Range("D7").Value = 62.17

# ---- New column F:I headers (mirrors A:D row 1) ----
$ws.Range("F1").Value = "5m"
$ws.Range("G1").Value = "30m"
$ws.Range("H1").Value = "1H"
$ws.Range("I1").Value = "4H"

# Row 2 label
$ws.Range("F2").Value = "Khoi luong"
$ws.Range("G2").Value = "Khoi luong"
$ws.Range("H2").Value = "Khoi luong"
$ws.Range("I2").Value = "Khoi luong"

# Row 3 values
$ws.Range("F3").Value = 371
$ws.Range("G3").Value = 857
$ws.Range("H3").Value = 857
$ws.Range("I3").Value = 837

# Row 4 values
$ws.Range("F4").Value = 472.85
$ws.Range("G4").Value = 1926.35
$ws.Range("H4").Value = 4554.8
$ws.Range("I4").Value = 16735.849999999999

# Row 5 label
$ws.Range("F5").Value = "RSI"
$ws.Range("G5").Value = "RSI"
$ws.Range("H5").Value = "RSI"
$ws.Range("I5").Value = "RSI"

# Row 6 values
$ws.Range("F6").Value = 55.72
$ws.Range("G6").Value = 63.39
$ws.Range("H6").Value = 69.52
$ws.Range("I6").Value = 67.97

# Row 7 values
$ws.Range("F7").Value = 62.17
$ws.Range("G7").Value = 64.87
$ws.Range("H7").Value = 68.7
$ws.Range("I7").Value = 62.25

# Row 8 label (empty set symbol)
$ws.Range("F8").Value = [char]0x2205
$ws.Range("G8").Value = [char]0x2205
$ws.Range("H8").Value = [char]0x2205
$ws.Range("I8").Value = [char]0x2205

# Row 9 label (empty set symbol)
$ws.Range("F9").Value = [char]0x2205
$ws.Range("G9").Value = [char]0x2205
$ws.Range("H9").Value = [char]0x2205
$ws.Range("I9").Value = [char]0x2205

# ---- Column widths for F:I (closest achievable values; runtime stores
#      width on a coarser grid than native Excel's 1/256-char units) ----
$ws.Columns.Item(6).ColumnWidth = 9.65
$ws.Columns.Item(7).ColumnWidth = 10.17
$ws.Columns.Item(8).ColumnWidth = 9.45
$ws.Columns.Item(9).ColumnWidth = 9.65

# ---- New formulas G39, G40, G41, G42 ----
$ws.Range("G39").Formula = "=A4-F4"
$ws.Range("G40").Formula = "=B4-G4"
$ws.Range("G41").Formula = "=C4-H4"
$ws.Range("G42").Formula = "=D4-I4"

# ---- View state: scroll + selection ----
$excel.ActiveWindow.ScrollRow = 29
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("G43").Select()
